$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object 'object[,]' 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.044223236010621
$dataBF[0,2] = 1.051257324750364
$dataBF[0,3] = 1.057403316258272
$dataBF[0,4] = 1.063932650317414
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.045060905255169
$dataBF[1,2] = 1.051914999137675
$dataBF[1,3] = 1.058202620256464
$dataBF[1,4] = 1.06475528482229
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.04560356313494
$dataBF[2,2] = 1.052341061988759
$dataBF[2,3] = 1.058720786999058
$dataBF[2,4] = 1.065288573684773
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.045831845623656
$dataBF[3,2] = 1.052520297563154
$dataBF[3,3] = 1.058938853072243
$dataBF[3,4] = 1.065513003201093
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.045870183953602
$dataBF[4,2] = 1.052550398916118
$dataBF[4,3] = 1.058975480685226
$dataBF[4,4] = 1.065550699620548
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.045606612872896
$dataBF[5,2] = 1.052343456478924
$dataBF[5,3] = 1.058723699912077
$dataBF[5,4] = 1.065291571602516
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.044506198828215
$dataBF[6,2] = 1.051479483546067
$dataBF[6,3] = 1.057673244114245
$dataBF[6,4] = 1.064210457468981
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.042572026839891
$dataBF[7,2] = 1.049960988851936
$dataBF[7,3] = 1.055829671544742
$dataBF[7,4] = 1.062313060384124
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.041285971865566
$dataBF[8,2] = 1.048951412619736
$dataBF[8,3] = 1.054605756293536
$dataBF[8,4] = 1.06105339940752
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.040729921624822
$dataBF[9,2] = 1.048514930730263
$dataBF[9,3] = 1.054077029120304
$dataBF[9,4] = 1.060509227185009
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.04052350469391
$dataBF[10,2] = 1.048352904609247
$dataBF[10,3] = 1.05388082398899
$dataBF[10,4] = 1.060307290100747
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.040567776147847
$dataBF[11,2] = 1.048387655085977
$dataBF[11,3] = 1.053922902131269
$dataBF[11,4] = 1.060350597538064
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.040712856579523
$dataBF[12,2] = 1.048501535495167
$dataBF[12,3] = 1.054060806902968
$dataBF[12,4] = 1.060492531059804
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.040802262020836
$dataBF[13,2] = 1.048571714635671
$dataBF[13,3] = 1.054145799512479
$dataBF[13,4] = 1.060580006582836
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.041322892509873
$dataBF[14,2] = 1.048980394785007
$dataBF[14,3] = 1.054640872387177
$dataBF[14,4] = 1.061089541240012
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.041649690998856
$dataBF[15,2] = 1.049236930062452
$dataBF[15,3] = 1.054951750870911
$dataBF[15,4] = 1.061409500114817
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.041840386038354
$dataBF[16,2] = 1.049386627395553
$dataBF[16,3] = 1.055133200175761
$dataBF[16,4] = 1.06159624904673
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.041905421469092
$dataBF[17,2] = 1.049437681226846
$dataBF[17,3] = 1.05519508984916
$dataBF[17,4] = 1.061659946295161
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.041614620429269
$dataBF[18,2] = 1.049209399555131
$dataBF[18,3] = 1.054918384213897
$dataBF[18,4] = 1.061375158887569
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.040670130547414
$dataBF[19,2] = 1.048467997703668
$dataBF[19,3] = 1.054020192188232
$dataBF[19,4] = 1.060450729839407
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.040077015663062
$dataBF[20,2] = 1.048002443955778
$dataBF[20,3] = 1.053456550052606
$dataBF[20,4] = 1.05987062055295
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.040391367953339
$dataBF[21,2] = 1.048249185792867
$dataBF[21,3] = 1.053755243740353
$dataBF[21,4] = 1.060178041011479
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.041630467066147
$dataBF[22,2] = 1.049221839206845
$dataBF[22,3] = 1.054933460802065
$dataBF[22,4] = 1.06139067583223
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.043071465933675
$dataBF[23,2] = 1.050353078636142
$dataBF[23,3] = 1.056305382856304
$dataBF[23,4] = 1.062802662662131

$ws.Range("B2:F25").Value = $dataBF

$dataIN = New-Object 'object[,]' 24,6
$dataIN[0,0] = 1.045411574056603
$dataIN[0,1] = 1.049289472796228
$dataIN[0,2] = 1.054009328123815
$dataIN[0,3] = 1.060138381339192
$dataIN[0,4] = 1.066649952673669
$dataIN[0,5] = 1.020308718565064
$dataIN[1,0] = 1.045628548531108
$dataIN[1,1] = 1.049774863374952
$dataIN[1,2] = 1.054480011405298
$dataIN[1,3] = 1.060751551630359
$dataIN[1,4] = 1.067287675183049
$dataIN[1,5] = 1.020470953976256
$dataIN[2,0] = 1.045767903474479
$dataIN[2,1] = 1.050088892953594
$dataIN[2,2] = 1.054784392292347
$dataIN[2,3] = 1.061148619222951
$dataIN[2,4] = 1.06770065925739
$dataIN[2,5] = 1.020575881784666
$dataIN[3,0] = 1.045826238001789
$dataIN[3,1] = 1.050220897592088
$dataIN[3,2] = 1.054912308950255
$dataIN[3,3] = 1.061315617993993
$dataIN[3,4] = 1.067874356319805
$dataIN[3,5] = 1.020619981148374
$dataIN[4,0] = 1.045836017924563
$dataIN[4,1] = 1.050243060928763
$dataIN[4,2] = 1.054933784038394
$dataIN[4,3] = 1.061343661973108
$dataIN[4,4] = 1.067903525380017
$dataIN[4,5] = 1.020627384891339
$dataIN[5,0] = 1.045768683927785
$dataIN[5,1] = 1.050090656860286
$dataIN[5,2] = 1.054786101699748
$dataIN[5,3] = 1.0611508503895
$dataIN[5,4] = 1.067702979899413
$dataIN[5,5] = 1.020576471091394
$dataIN[6,0] = 1.045485116864976
$dataIN[6,1] = 1.049453522514809
$dataIN[6,2] = 1.054168434762643
$dataIN[6,3] = 1.060345540933613
$dataIN[6,4] = 1.066865403641789
$dataIN[6,5] = 1.020363556716493
$dataIN[7,0] = 1.044977490991519
$dataIN[7,1] = 1.048330479883007
$dataIN[7,2] = 1.053078686712963
$dataIN[7,3] = 1.058928890910801
$dataIN[7,4] = 1.065392124974479
$dataIN[7,5] = 1.01998801728361
$dataIN[8,0] = 1.044633781345611
$dataIN[8,1] = 1.047581639265308
$dataIN[8,2] = 1.052351366688859
$dataIN[8,3] = 1.057986169748655
$dataIN[8,4] = 1.064411810107189
$dataIN[8,5] = 1.019737445706467
$dataIN[9,0] = 1.044483706300983
$dataIN[9,1] = 1.0472573646577
$dataIN[9,2] = 1.052036251070023
$dataIN[9,3] = 1.057578386829261
$dataIN[9,4] = 1.063987787271462
$dataIN[9,5] = 1.019628901109619
$dataIN[10,0] = 1.044427775288252
$dataIN[10,1] = 1.047136912671292
$dataIN[10,2] = 1.051919177216217
$dataIN[10,3] = 1.057426983016059
$dataIN[10,4] = 1.063830357075783
$dataIN[10,5] = 1.019588576462165
$dataIN[11,0] = 1.044439781104934
$dataIN[11,1] = 1.047162750099057
$dataIN[11,2] = 1.051944291104868
$dataIN[11,3] = 1.057459456684762
$dataIN[11,4] = 1.063864123158252
$dataIN[11,5] = 1.019597226515679
$dataIN[12,0] = 1.044479086825019
$dataIN[12,1] = 1.047247408091379
$dataIN[12,2] = 1.052026574228943
$dataIN[12,3] = 1.057565870409329
$dataIN[12,4] = 1.063974772589191
$dataIN[12,5] = 1.019625567988486
$dataIN[13,0] = 1.044503279693143
$dataIN[13,1] = 1.047299568453784
$dataIN[13,2] = 1.052077268188018
$dataIN[13,3] = 1.057631444066989
$dataIN[13,4] = 1.064042956782188
$dataIN[13,5] = 1.019643029277936
$dataIN[14,0] = 1.044643715139141
$dataIN[14,1] = 1.047603159957106
$dataIN[14,2] = 1.05237227617105
$dataIN[14,3] = 1.05801324197514
$dataIN[14,4] = 1.064439960939052
$dataIN[14,5] = 1.019744648537028
$dataIN[15,0] = 1.04473147343448
$dataIN[15,1] = 1.047793589998773
$dataIN[15,2] = 1.052557279311121
$dataIN[15,3] = 1.058252847699325
$dataIN[15,4] = 1.064689115477009
$dataIN[15,5] = 1.01980837978211
$dataIN[16,0] = 1.044782541039601
$dataIN[16,1] = 1.047904662441532
$dataIN[16,2] = 1.052665170934442
$dataIN[16,3] = 1.058392646233105
$dataIN[16,4] = 1.064834487372885
$dataIN[16,5] = 1.019845548721832
$dataIN[17,0] = 1.044799933339515
$dataIN[17,1] = 1.047942534865101
$dataIN[17,2] = 1.052701956167226
$dataIN[17,3] = 1.058440320753353
$dataIN[17,4] = 1.064884062890964
$dataIN[17,5] = 1.019858221605974
$dataIN[18,0] = 1.044722070239713
$dataIN[18,1] = 1.047773158873503
$dataIN[18,2] = 1.05253743203352
$dataIN[18,3] = 1.058227136065087
$dataIN[18,4] = 1.06466237896211
$dataIN[18,5] = 1.019801542476117
$dataIN[19,0] = 1.044467517410766
$dataIN[19,1] = 1.047222478466333
$dataIN[19,2] = 1.052002344603413
$dataIN[19,3] = 1.057534532414918
$dataIN[19,4] = 1.063942187122855
$dataIN[19,5] = 1.019617222302632
$dataIN[20,0] = 1.044306391553266
$dataIN[20,1] = 1.046876233147457
$dataIN[20,2] = 1.051665764817845
$dataIN[20,3] = 1.057099441399392
$dataIN[20,4] = 1.063489784175582
$dataIN[20,5] = 1.019501296383569
$dataIN[21,0] = 1.044391909365855
$dataIN[21,1] = 1.047059784922009
$dataIN[21,2] = 1.051844205811505
$dataIN[21,3] = 1.057330055172713
$dataIN[21,4] = 1.063729572096993
$dataIN[21,5] = 1.019562754208171
$dataIN[22,0] = 1.044726319510551
$dataIN[22,1] = 1.047782390826459
$dataIN[22,2] = 1.052546400218435
$dataIN[22,3] = 1.058238753920705
$dataIN[22,4] = 1.064674459905759
$dataIN[22,5] = 1.019804631974096
$dataIN[23,0] = 1.04510966086373
$dataIN[23,1] = 1.048620844588183
$dataIN[23,2] = 1.053360563296764
$dataIN[23,3] = 1.059294833928986
$dataIN[23,4] = 1.06577268019763
$dataIN[23,5] = 1.020085142646241

$ws.Range("I2:N25").Value = $dataIN

Write-Host "done"